$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the two runs that make up "Fri Sep 14" / " 11:44:22 PDT 2017"
#    into a single run with the full text "Fri Sep 14 11:44:22 PDT 2017".
#    A same-text Find/Replace over the whole (already correct) text
#    causes the engine to re-serialize it as one run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Fri Sep 14 11:44:22 PDT 2017", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Fri Sep 14 11:44:22 PDT 2017", 2)

# ---------------------------------------------------------------------
# 2) Append a brand-new purchase record (dated "Sun Sep 16 ...") right
#    after the last existing record ("Amount Received mode - CASH"),
#    before the trailing blank paragraphs that close out the document.
# ---------------------------------------------------------------------

# Locate the paragraph that holds the last "Amount Received mode ... - CASH" line.
$targetIndex = -1
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -match "Amount Received mode") {
        $targetIndex = $i
        break
    }
}

$lines = @(
    @{ Text = ""; Bold = $false },
    @{ Text = "Sun Sep 16 11:33:19 PDT 2017"; Bold = $false },
    @{ Text = "Person Name`t`t`t`t- HSJ"; Bold = $false },
    @{ Text = "---------------------------------------------------------------"; Bold = $false },
    @{ Text = "Item Name`t`t`t`t- CARROT EVE"; Bold = $false },
    @{ Text = "Number of Pockets`t`t`t- 1"; Bold = $false },
    @{ Text = "Number of KGs`t`t`t- 97"; Bold = $false },
    @{ Text = "Rate`t`t`t`t- 20"; Bold = $false },
    @{ Text = "Total Price`t`t`t`t- 1940.0"; Bold = $false },
    @{ Text = "Amount balance`t`t`t- 18664.0"; Bold = $true },
    @{ Text = ""; Bold = $false }
)

$currentIndex = $targetIndex
foreach ($line in $lines) {
    $currentPara = $d.Paragraphs($currentIndex)
    $currentPara.Range.InsertParagraphAfter()
    $currentIndex = $currentIndex + 1
    $newPara = $d.Paragraphs($currentIndex)
    if ($line.Text.Length -gt 0) {
        $newPara.Range.Text = $line.Text
    }
    if ($line.Bold) {
        $newPara.Range.Font.Bold = 1
    }
}

Write-Output "done"
